# OCPTIME_ASSEMBLE_MAT_FOR_LS -> OCPTIME_CONVERT_MAT_FOR_LS_IF
# Scale down the "N" (number of blocks) columns by a factor of 3 (dividing
# dx/dy counts) while keeping per-block sizes but shrinking the overall grid
# size (Dx, Dy, Dz) accordingly, and fix up the dz-derived column which isn't
# an exact integer ratio.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Grid extents (B2:D2) ---
$ws.Range("B2").Value = 10000
$ws.Range("C2").Value = 10000
$ws.Range("D2").Value = 100

# --- Row 2 counts (L2, N2, P2) ---
$ws.Range("L2").Value = 1
$ws.Range("N2").Value = 1
$ws.Range("P2").Value = 1

# --- Row 3 counts (L3, N3, P3) ---
$ws.Range("L3").Value = 2
$ws.Range("N3").Value = 2
$ws.Range("P3").Value = 2

# --- Row 4 counts (L4, N4, P4) ---
$ws.Range("L4").Value = 4
$ws.Range("N4").Value = 4
$ws.Range("P4").Value = 4

# --- Row 5 counts (L5, N5, P5) ---
$ws.Range("L5").Value = 8
$ws.Range("N5").Value = 8
$ws.Range("P5").Value = 8.3

# --- Row 6 counts (L6, N6, P6) ---
$ws.Range("L6").Value = 16
$ws.Range("N6").Value = 16
$ws.Range("P6").Value = 16.7

# --- Selection moves from K11 to M11 ---
$ws.Range("M11").Select()

# --- Window height grows slightly ---
$excel.ActiveWindow.Height = 9925
